$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test18042025@289.com"
$ws.Range("B2").Value = "Test18042025@289"
